$wb = $excel.ActiveWorkbook

# --- Step 1: rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 2: add a new "PO Forecast" sheet as the 3rd (last) sheet ---
# Duplicate the "Weekly Quantity" sheet so the new sheet inherits the same
# sheet-level properties (outline/page-setup, page margins) and the same
# header / date-column cell styles already used in this workbook, then wipe
# its values and rebuild them from scratch.
$wsWeekly.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.ClearContents()

# Extend the header style (bold / border / centered) to the new C & D columns,
# and the date-serial style to the full A2:A75 data range.
$wsForecast.Range("A1:B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)
$wsForecast.Range("A2").Copy()
$wsForecast.Range("A2:A75").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Step 4: data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$wsForecast.Range("A2").Value = 44955.99999999999
$wsForecast.Range("B2").Value = 64
$wsForecast.Range("C2").Value = -62.5817767187154
$wsForecast.Range("D2").Value = 206.6564482053743
$wsForecast.Range("A3").Value = 44969.99999999999
$wsForecast.Range("B3").Value = 66
$wsForecast.Range("C3").Value = -61.51296061054284
$wsForecast.Range("D3").Value = 193.8596457573256
$wsForecast.Range("A4").Value = 44976.99999999999
$wsForecast.Range("B4").Value = 67
$wsForecast.Range("C4").Value = -68.92147015509397
$wsForecast.Range("D4").Value = 197.6974743370473
$wsForecast.Range("A5").Value = 44983.99999999999
$wsForecast.Range("B5").Value = 68
$wsForecast.Range("C5").Value = -60.14984545463125
$wsForecast.Range("D5").Value = 192.089848385362
$wsForecast.Range("A6").Value = 44990.99999999999
$wsForecast.Range("B6").Value = 69
$wsForecast.Range("C6").Value = -77.99123027054029
$wsForecast.Range("D6").Value = 188.4604455468331
$wsForecast.Range("A7").Value = 44997.99999999999
$wsForecast.Range("B7").Value = 70
$wsForecast.Range("C7").Value = -67.50259711274333
$wsForecast.Range("D7").Value = 191.9130157968731
$wsForecast.Range("A8").Value = 45004.99999999999
$wsForecast.Range("B8").Value = 71
$wsForecast.Range("C8").Value = -57.33786344199211
$wsForecast.Range("D8").Value = 202.3923485696238
$wsForecast.Range("A9").Value = 45011.99999999999
$wsForecast.Range("B9").Value = 72
$wsForecast.Range("C9").Value = -54.39284248652898
$wsForecast.Range("D9").Value = 202.3910036591091
$wsForecast.Range("A10").Value = 45025.99999999999
$wsForecast.Range("B10").Value = 74
$wsForecast.Range("C10").Value = -64.65652752323076
$wsForecast.Range("D10").Value = 197.6638626324447
$wsForecast.Range("A11").Value = 45039.99999999999
$wsForecast.Range("B11").Value = 77
$wsForecast.Range("C11").Value = -54.59105375392785
$wsForecast.Range("D11").Value = 210.073934578274
$wsForecast.Range("A12").Value = 45081.99999999999
$wsForecast.Range("B12").Value = 83
$wsForecast.Range("C12").Value = -40.32358154868776
$wsForecast.Range("D12").Value = 215.5470470089504
$wsForecast.Range("A13").Value = 45088.99999999999
$wsForecast.Range("B13").Value = 84
$wsForecast.Range("C13").Value = -52.24160808383002
$wsForecast.Range("D13").Value = 208.9631609375414
$wsForecast.Range("A14").Value = 45095.99999999999
$wsForecast.Range("B14").Value = 85
$wsForecast.Range("C14").Value = -48.96925758907283
$wsForecast.Range("D14").Value = 213.729479888269
$wsForecast.Range("A15").Value = 45102.99999999999
$wsForecast.Range("B15").Value = 86
$wsForecast.Range("C15").Value = -56.97617773807804
$wsForecast.Range("D15").Value = 216.3396249962008
$wsForecast.Range("A16").Value = 45109.99999999999
$wsForecast.Range("B16").Value = 87
$wsForecast.Range("C16").Value = -41.19952644774843
$wsForecast.Range("D16").Value = 216.9264697787791
$wsForecast.Range("A17").Value = 45130.99999999999
$wsForecast.Range("B17").Value = 91
$wsForecast.Range("C17").Value = -37.72308230005661
$wsForecast.Range("D17").Value = 215.7441311545828
$wsForecast.Range("A18").Value = 45137.99999999999
$wsForecast.Range("B18").Value = 92
$wsForecast.Range("C18").Value = -47.51318435767615
$wsForecast.Range("D18").Value = 230.9679779078472
$wsForecast.Range("A19").Value = 45144.99999999999
$wsForecast.Range("B19").Value = 93
$wsForecast.Range("C19").Value = -40.13887235000918
$wsForecast.Range("D19").Value = 232.8451752806505
$wsForecast.Range("A20").Value = 45151.99999999999
$wsForecast.Range("B20").Value = 94
$wsForecast.Range("C20").Value = -40.4039010060068
$wsForecast.Range("D20").Value = 220.9533161616009
$wsForecast.Range("A21").Value = 45158.99999999999
$wsForecast.Range("B21").Value = 95
$wsForecast.Range("C21").Value = -32.63243398782489
$wsForecast.Range("D21").Value = 224.6264695598984
$wsForecast.Range("A22").Value = 45165.99999999999
$wsForecast.Range("B22").Value = 96
$wsForecast.Range("C22").Value = -36.0707487031124
$wsForecast.Range("D22").Value = 222.5690154624815
$wsForecast.Range("A23").Value = 45172.99999999999
$wsForecast.Range("B23").Value = 97
$wsForecast.Range("C23").Value = -22.44346996032181
$wsForecast.Range("D23").Value = 225.2607723715232
$wsForecast.Range("A24").Value = 45179.99999999999
$wsForecast.Range("B24").Value = 98
$wsForecast.Range("C24").Value = -19.24924286980228
$wsForecast.Range("D24").Value = 220.3669700017284
$wsForecast.Range("A25").Value = 45186.99999999999
$wsForecast.Range("B25").Value = 99
$wsForecast.Range("C25").Value = -23.9713111144198
$wsForecast.Range("D25").Value = 227.7050102010936
$wsForecast.Range("A26").Value = 45193.99999999999
$wsForecast.Range("B26").Value = 101
$wsForecast.Range("C26").Value = -30.19750222649008
$wsForecast.Range("D26").Value = 224.1340876276751
$wsForecast.Range("A27").Value = 45200.99999999999
$wsForecast.Range("B27").Value = 102
$wsForecast.Range("C27").Value = -33.02058657327188
$wsForecast.Range("D27").Value = 228.4342851569741
$wsForecast.Range("A28").Value = 45207.99999999999
$wsForecast.Range("B28").Value = 103
$wsForecast.Range("C28").Value = -23.37756932229329
$wsForecast.Range("D28").Value = 238.420362940408
$wsForecast.Range("A29").Value = 45214.99999999999
$wsForecast.Range("B29").Value = 104
$wsForecast.Range("C29").Value = -26.65895609714185
$wsForecast.Range("D29").Value = 226.8826269770099
$wsForecast.Range("A30").Value = 45221.99999999999
$wsForecast.Range("B30").Value = 105
$wsForecast.Range("C30").Value = -21.75068446079355
$wsForecast.Range("D30").Value = 244.2779173843764
$wsForecast.Range("A31").Value = 45242.99999999999
$wsForecast.Range("B31").Value = 108
$wsForecast.Range("C31").Value = -24.69336863306893
$wsForecast.Range("D31").Value = 239.4898994211561
$wsForecast.Range("A32").Value = 45249.99999999999
$wsForecast.Range("B32").Value = 109
$wsForecast.Range("C32").Value = -26.40385448460479
$wsForecast.Range("D32").Value = 234.3654733103801
$wsForecast.Range("A33").Value = 45256.99999999999
$wsForecast.Range("B33").Value = 110
$wsForecast.Range("C33").Value = -21.33557879771219
$wsForecast.Range("D33").Value = 236.0188142917754
$wsForecast.Range("A34").Value = 45263.99999999999
$wsForecast.Range("B34").Value = 111
$wsForecast.Range("C34").Value = -19.14321081835924
$wsForecast.Range("D34").Value = 253.4857303210721
$wsForecast.Range("A35").Value = 45270.99999999999
$wsForecast.Range("B35").Value = 112
$wsForecast.Range("C35").Value = -15.81905233892271
$wsForecast.Range("D35").Value = 251.4335631236878
$wsForecast.Range("A36").Value = 45298.99999999999
$wsForecast.Range("B36").Value = 117
$wsForecast.Range("C36").Value = -18.34890084317227
$wsForecast.Range("D36").Value = 246.6517600173265
$wsForecast.Range("A37").Value = 45305.99999999999
$wsForecast.Range("B37").Value = 118
$wsForecast.Range("C37").Value = -3.616699857979824
$wsForecast.Range("D37").Value = 242.5153703804851
$wsForecast.Range("A38").Value = 45312.99999999999
$wsForecast.Range("B38").Value = 119
$wsForecast.Range("C38").Value = -15.80305276675174
$wsForecast.Range("D38").Value = 253.8019887590873
$wsForecast.Range("A39").Value = 45319.99999999999
$wsForecast.Range("B39").Value = 120
$wsForecast.Range("C39").Value = -13.4468107510624
$wsForecast.Range("D39").Value = 259.0424712298944
$wsForecast.Range("A40").Value = 45326.99999999999
$wsForecast.Range("B40").Value = 121
$wsForecast.Range("C40").Value = -14.65908228496219
$wsForecast.Range("D40").Value = 253.1509140886062
$wsForecast.Range("A41").Value = 45333.99999999999
$wsForecast.Range("B41").Value = 122
$wsForecast.Range("C41").Value = -17.20925053349944
$wsForecast.Range("D41").Value = 250.6443547619062
$wsForecast.Range("A42").Value = 45347.99999999999
$wsForecast.Range("B42").Value = 124
$wsForecast.Range("C42").Value = -14.90125631143898
$wsForecast.Range("D42").Value = 251.1831797684817
$wsForecast.Range("A43").Value = 45354.99999999999
$wsForecast.Range("B43").Value = 125
$wsForecast.Range("C43").Value = -3.443295766636244
$wsForecast.Range("D43").Value = 258.2258794738868
$wsForecast.Range("A44").Value = 45361.99999999999
$wsForecast.Range("B44").Value = 127
$wsForecast.Range("C44").Value = -3.707204883791736
$wsForecast.Range("D44").Value = 260.2610301324789
$wsForecast.Range("A45").Value = 45368.99999999999
$wsForecast.Range("B45").Value = 128
$wsForecast.Range("C45").Value = 1.339416512816808
$wsForecast.Range("D45").Value = 255.2739192339698
$wsForecast.Range("A46").Value = 45375.99999999999
$wsForecast.Range("B46").Value = 129
$wsForecast.Range("C46").Value = 0.7527636106177328
$wsForecast.Range("D46").Value = 259.1985820254783
$wsForecast.Range("A47").Value = 45382.99999999999
$wsForecast.Range("B47").Value = 130
$wsForecast.Range("C47").Value = -7.57734573771845
$wsForecast.Range("D47").Value = 251.4600586152889
$wsForecast.Range("A48").Value = 45389.99999999999
$wsForecast.Range("B48").Value = 131
$wsForecast.Range("C48").Value = 5.163476627552382
$wsForecast.Range("D48").Value = 262.7333779632535
$wsForecast.Range("A49").Value = 45396.99999999999
$wsForecast.Range("B49").Value = 132
$wsForecast.Range("C49").Value = 16.97018830979288
$wsForecast.Range("D49").Value = 264.617671982717
$wsForecast.Range("A50").Value = 45403.99999999999
$wsForecast.Range("B50").Value = 133
$wsForecast.Range("C50").Value = -1.944800912658563
$wsForecast.Range("D50").Value = 259.0597491265291
$wsForecast.Range("A51").Value = 45410.99999999999
$wsForecast.Range("B51").Value = 134
$wsForecast.Range("C51").Value = 6.09596938574371
$wsForecast.Range("D51").Value = 258.7392067987943
$wsForecast.Range("A52").Value = 45417.99999999999
$wsForecast.Range("B52").Value = 135
$wsForecast.Range("C52").Value = -0.2229395088558367
$wsForecast.Range("D52").Value = 258.6068289497743
$wsForecast.Range("A53").Value = 45424.99999999999
$wsForecast.Range("B53").Value = 136
$wsForecast.Range("C53").Value = -1.875506355958767
$wsForecast.Range("D53").Value = 274.7631302700093
$wsForecast.Range("A54").Value = 45431.99999999999
$wsForecast.Range("B54").Value = 137
$wsForecast.Range("C54").Value = 10.96970428460158
$wsForecast.Range("D54").Value = 281.5953167231898
$wsForecast.Range("A55").Value = 45438.99999999999
$wsForecast.Range("B55").Value = 139
$wsForecast.Range("C55").Value = 8.21558082468027
$wsForecast.Range("D55").Value = 275.2117814632168
$wsForecast.Range("A56").Value = 45445.99999999999
$wsForecast.Range("B56").Value = 140
$wsForecast.Range("C56").Value = 6.330062896631408
$wsForecast.Range("D56").Value = 261.4029578017483
$wsForecast.Range("A57").Value = 45459.99999999999
$wsForecast.Range("B57").Value = 142
$wsForecast.Range("C57").Value = 5.611068825818452
$wsForecast.Range("D57").Value = 266.3178106123411
$wsForecast.Range("A58").Value = 45466.99999999999
$wsForecast.Range("B58").Value = 143
$wsForecast.Range("C58").Value = 19.462712598406
$wsForecast.Range("D58").Value = 264.2110147452472
$wsForecast.Range("A59").Value = 45487.99999999999
$wsForecast.Range("B59").Value = 146
$wsForecast.Range("C59").Value = 12.56009387886268
$wsForecast.Range("D59").Value = 272.0888013843039
$wsForecast.Range("A60").Value = 45494.99999999999
$wsForecast.Range("B60").Value = 147
$wsForecast.Range("C60").Value = 18.61580252772221
$wsForecast.Range("D60").Value = 274.8346712429371
$wsForecast.Range("A61").Value = 45529.99999999999
$wsForecast.Range("B61").Value = 153
$wsForecast.Range("C61").Value = 25.63413166597774
$wsForecast.Range("D61").Value = 285.1044703297391
$wsForecast.Range("A62").Value = 45536.99999999999
$wsForecast.Range("B62").Value = 154
$wsForecast.Range("C62").Value = 24.10291535891425
$wsForecast.Range("D62").Value = 285.1428769699596
$wsForecast.Range("A63").Value = 45543.99999999999
$wsForecast.Range("B63").Value = 155
$wsForecast.Range("C63").Value = 26.15445774146766
$wsForecast.Range("D63").Value = 288.9064204156553
$wsForecast.Range("A64").Value = 45578.99999999999
$wsForecast.Range("B64").Value = 160
$wsForecast.Range("C64").Value = 31.51452335486405
$wsForecast.Range("D64").Value = 293.3888626660948
$wsForecast.Range("A65").Value = 45585.99999999999
$wsForecast.Range("B65").Value = 161
$wsForecast.Range("C65").Value = 28.97100748780574
$wsForecast.Range("D65").Value = 290.1690277030777
$wsForecast.Range("A66").Value = 45592.99999999999
$wsForecast.Range("B66").Value = 162
$wsForecast.Range("C66").Value = 38.54276402235418
$wsForecast.Range("D66").Value = 293.4041934023787
$wsForecast.Range("A67").Value = 45599.99999999999
$wsForecast.Range("B67").Value = 164
$wsForecast.Range("C67").Value = 40.83724675422823
$wsForecast.Range("D67").Value = 295.0317616795652
$wsForecast.Range("A68").Value = 45606.99999999999
$wsForecast.Range("B68").Value = 165
$wsForecast.Range("C68").Value = 29.19154194450919
$wsForecast.Range("D68").Value = 297.3489442462069
$wsForecast.Range("A69").Value = 45613.99999999999
$wsForecast.Range("B69").Value = 166
$wsForecast.Range("C69").Value = 32.62367560475913
$wsForecast.Range("D69").Value = 292.2728568999059
$wsForecast.Range("A70").Value = 45620.99999999999
$wsForecast.Range("B70").Value = 167
$wsForecast.Range("C70").Value = 40.39930943829091
$wsForecast.Range("D70").Value = 297.1018638577021
$wsForecast.Range("A71").Value = 45627.99999999999
$wsForecast.Range("B71").Value = 168
$wsForecast.Range("C71").Value = 28.36172672892165
$wsForecast.Range("D71").Value = 296.1096836276216
$wsForecast.Range("A72").Value = 45634.99999999999
$wsForecast.Range("B72").Value = 169
$wsForecast.Range("C72").Value = 39.93747989606688
$wsForecast.Range("D72").Value = 298.7302370188145
$wsForecast.Range("A73").Value = 45641.99999999999
$wsForecast.Range("B73").Value = 170
$wsForecast.Range("C73").Value = 44.51232599757079
$wsForecast.Range("D73").Value = 301.3165623926762
$wsForecast.Range("A74").Value = 45648.99999999999
$wsForecast.Range("B74").Value = 171
$wsForecast.Range("C74").Value = 39.85373666743065
$wsForecast.Range("D74").Value = 304.5956912629257
$wsForecast.Range("A75").Value = 45655.99999999999
$wsForecast.Range("B75").Value = 172
$wsForecast.Range("C75").Value = 41.01412364045583
$wsForecast.Range("D75").Value = 302.6392588065156

Write-Host "PO Forecast sheet created with" ($wsForecast.UsedRange.Rows.Count) "rows and" ($wsForecast.UsedRange.Columns.Count) "cols"
